# Apply updated crypto price/volume data (and two coin-row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.321.93'
$ws.Range('E2').Value = '  -2.20%  '
$ws.Range('D3').Value = '2.639.07'
$ws.Range('E3').Value = '  -3.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '548.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.67%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -1.80%  '
$ws.Range('E9').Value = '  -4.49%  '
$ws.Range('E10').Value = '  -4.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.42'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.46%  '
$ws.Range('E12').Value = '  -4.70%  '
$ws.Range('D13').Value = '3.104.47'
$ws.Range('E14').Value = '  -4.99%  '
$ws.Range('D15').Value = '62.235.16'
$ws.Range('E15').Value = '  -2.09%  '
$ws.Range('E16').Value = '  -4.04%  '
$ws.Range('D17').Value = '2.639.72'
$ws.Range('E17').Value = '  -3.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.65'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '338.39'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.54%  '
$ws.Range('E21').Value = '  -8.10%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.498'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.59'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.168'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  -4.34%  '
$ws.Range('D28').Value = '0.0₃0832'
$ws.Range('E28').Value = '  -8.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('E30').Value = '  -1.35%  '
$ws.Range('E31').Value = '  -5.45%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '160.16'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.98%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.72'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.86%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '19.18'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.45%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.42'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.07%  '
$ws.Range('E37').Value = '  -4.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '332.43'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.07'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.93%  '
$ws.Range('E40').Value = '  -7.41%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.92'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.95%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.90'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.71%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.607'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.60'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0544'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.67%  '
$ws.Range('E49').Value = '  -3.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '127.15'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0236'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.41%  '
